$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 179, shifting existing rows 179:210 down to 180:211
$ws.Rows.Item(179).Insert()

# Populate the newly inserted row 179 with the new record's data
$ws.Range("A179").Value = 5
$ws.Range("B179").Value = "Macroferia Regional de Talca"
$ws.Range("C179").Value = "Maule"
$ws.Range("D179").Value = 44522
$ws.Range("E179").Value = 7
$ws.Range("F179").Value = 100114014
$ws.Range("G179").Value = "Betarraga"
$ws.Range("H179").Value = "Sin especificar"
$ws.Range("I179").Value = "Primera"
$ws.Range("J179").Value = 5000
$ws.Range("K179").Value = 500
$ws.Range("L179").Value = 500
$ws.Range("M179").Value = 500
$ws.Range("N179").Value = '$/paquete 5 unidades'
$ws.Range("O179").Value = "Región del Maule"
$ws.Range("P179").Value = 100
$ws.Range("Q179").Value = 5
$ws.Range("R179").Value = "Hortaliza"
